$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://fhir.nmdp.org/ig/matchsync/ValueSet/ms-abo-group-codes"
$meta.Range("B5").Value = "MatchSync ABO Group Value Set"
$meta.Range("B8").Value = "2023-01-12T09:36:27-06:00"
$meta.Range("B11").Value = "MatchSync ABO group codes. Combines LOINC, Snomed, and NMDP codes"

# --- Include ValueSets sheet (nmdp-abo-codes URL) ---
$vs1 = $wb.Worksheets.Item("Include ValueSets")
$vs1.Range("A2").Value = "http://fhir.nmdp.org/ig/matchsync/ValueSet/nmdp-abo-codes"

# --- Include ValueSets 3 sheet (sct-abo-group-codes URL) ---
$vs3 = $wb.Worksheets.Item("Include ValueSets 3")
$vs3.Range("A2").Value = "http://fhir.nmdp.org/ig/matchsync/ValueSet/sct-abo-group-codes"
